$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: alpha_distance_range -> Min 5 -> 5.4, Max 10.5 -> 10.199999999999999
$ws.Range("B2").Value = 5.4
$ws.Range("C2").Value = 10.199999999999999

# Row 3: beta_distance_range -> Min 5 -> 5.3 (Max 9 unchanged)
$ws.Range("B3").Value = 5.3

# Row 4: ratio_threshold_range -> Min 0.8 -> 0.85 (Max 1.25 unchanged)
$ws.Range("B4").Value = 0.85
